$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vocab words (column E, rows 13-23) - typed first, row by row,
# with rows 19/20 (slum/slump) filled out of order.
$ws.Range("E13").Value = "façade"
$ws.Range("E14").Value = "porch"
$ws.Range("E15").Value = "controversial"
$ws.Range("E16").Value = "derelict"
$ws.Range("E17").Value = "defunct"
$ws.Range("E18").Value = "estuary"
$ws.Range("E20").Value = "slump"
$ws.Range("E19").Value = "slum"
$ws.Range("E21").Value = "demote"
$ws.Range("E22").Value = "demobilise"
$ws.Range("E23").Value = "plead"

# New vocab words (column H, rows 5-7)
$ws.Range("H5").Value = "dilapidated"
$ws.Range("I5").Value = "(of a building or object) in a state of disrepair or ruin as a result of age or neglect."
$ws.Range("H6").Value = "amenity"
$ws.Range("H7").Value = "poise"

# Meanings for column E words (column F), filled after words, out of strict order.
$ws.Range("F14").Value = "cover area built onto house entrance"
$ws.Range("F15").Value = "causing a lot of disagreement/argument"
$ws.Range("F16").Value = "bad condition because of abandon/not used"
$ws.Range("F13").Value = "appearance | front of a building"
$ws.Range("F17").Value = "not working/existing now"
$ws.Range("F18").Value = "wide part of river before going into sea"
$ws.Range("F19").Value = "poor crowded area of a city"
$ws.Range("F20").Value = "fall down suddenly"
$ws.Range("F21").Value = "opposite of promote"
$ws.Range("F22").Value = "back off an army force"
$ws.Range("F23").Value = "legally appeal sth | ask sth"

# Meanings for column H words (column I, rows 6-7)
$ws.Range("I6").Value = "desireable/useful feature of a place"
$ws.Range("I7").Value = "behave in a calm/cofident way"

# Update the current selection to match the saved workbook state.
$ws.Range("I8").Select()
